# Update herolist to latest version: append newly released heroes.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("herolist")

$newHeroes = @(
    "npc_dota_hero_hoodwink",
    "npc_dota_hero_dawnbreaker",
    "npc_dota_hero_marci",
    "npc_dota_hero_primal_beast",
    "npc_dota_hero_muerta",
    "npc_dota_hero_ringmaster",
    "npc_dota_hero_kez"
)

$startRow = 122
for ($i = 0; $i -lt $newHeroes.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newHeroes[$i]
    $ws.Cells.Item($row, 2).Value = 1
}

# Scroll the view down and leave the selection on the last-added row, as in
# the authored workbook (topLeftCell="A85", activeCell="B124").
$excel.ActiveWindow.ScrollRow = 85
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B124").Select() | Out-Null
